# Regenerate merged AHB files
# - Rename header columns from *_old -> *_FV2304 and *_new -> *_FV2310
# - Wrap the data range in an Excel Table (ListObject)
# - Freeze the header row (pane split after row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "_old" headers (columns A-J) to "_FV2304"
for ($c = 1; $c -le 10; $c++) {
    $val = $ws.Cells.Item(1, $c).Value()
    if ($val -like "*_old") {
        $ws.Cells.Item(1, $c).Value = ($val -replace "_old$", "_FV2304")
    }
}

# 2) Rename the "_new" headers (columns L-U) to "_FV2310"
for ($c = 12; $c -le 21; $c++) {
    $val = $ws.Cells.Item(1, $c).Value()
    if ($val -like "*_new") {
        $ws.Cells.Item(1, $c).Value = ($val -replace "_new$", "_FV2310")
    }
}

# 3) Turn the used range into an Excel table
$rng = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 4) Freeze the top header row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$null = $excel.ActiveWindow.FreezePanes
